# SignalConfigs.xlsx update
# - Remove the "Testt" test row (was row 2) entirely.
# - Add three new test configs: TestAcc, TestDec, TimeTest as rows 6, 7, 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Testt" row (row 2) without shifting the rows below it up.
$ws.Rows(2).ClearContents()

# New row 6: TestAcc
$ws.Cells.Item(6, 1).Value = "TestAcc"
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(6, 3).Value = 10
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 20
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 1

# New row 7: TestDec
$ws.Cells.Item(7, 1).Value = "TestDec"
$ws.Cells.Item(7, 2).Value = 0
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = 10
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 20
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 1

# New row 8: TimeTest
$ws.Cells.Item(8, 1).Value = "TimeTest"
$ws.Cells.Item(8, 2).Value = 3
$ws.Cells.Item(8, 3).Value = 3
$ws.Cells.Item(8, 4).Value = 3
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 4
$ws.Cells.Item(8, 8).Value = 3
$ws.Cells.Item(8, 9).Value = 0.8

# Match the saved selection (J8) from the source workbook.
$ws.Range("J8").Select() | Out-Null
